$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtNum165 = "#,##0.##;""[""#,##0.##""]"";0"
$fmtNum2   = "0.00"

function Set-TextCell($ws, $addr, $value) {
    # Force the cell to store a shared-string (text) value instead of
    # auto-coercing numeric-looking strings ("1", "52.00", "2:0" ...) into
    # numbers. Setting NumberFormat to "@" immediately before the Value
    # assignment is what makes this engine keep it as text.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

function New-ProductRow($ws, $rowIndex, $srcRowIndex, $seq, $name, $stock, $order, $price, $sell, $deals) {
    # Insert a blank row, clone formatting (incl. merges) from a neighbouring
    # product row, then populate the cells.
    $ws.Rows.Item($rowIndex).Insert()

    $srcRange = "A" + $srcRowIndex + ":Q" + $srcRowIndex
    $dstRange = "A" + $rowIndex + ":Q" + $rowIndex
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    # PasteSpecial(formats) loses the readingOrder=RTL flag on column Q --
    # restore it explicitly.
    $addrQ = "Q" + $rowIndex
    $ws.Range($addrQ).ReadingOrder = 1

    $addrA = "A" + $rowIndex
    $ws.Range($addrA).Value = $seq

    $addrC = "C" + $rowIndex
    $addrH = "H" + $rowIndex
    $addrL = "L" + $rowIndex
    $addrN = "N" + $rowIndex
    $addrP = "P" + $rowIndex

    Set-TextCell $ws $addrC $name
    Set-TextCell $ws $addrH $stock
    Set-TextCell $ws $addrL $order
    Set-TextCell $ws $addrN $price
    Set-TextCell $ws $addrP $sell
    Set-TextCell $ws $addrQ $deals

    # Restore the per-column number formats that got clobbered by forcing "@".
    $ws.Range($addrL).NumberFormat = $fmtNum165
    $ws.Range($addrP).NumberFormat = $fmtNum2
}

# ---------------------------------------------------------------------------
# Step 1: insert TOBOLANZA before the "سرنجات 3 سم" row (currently row 14).
# Clone formatting from row 13 (T4-THYRO), the row immediately above.
# ---------------------------------------------------------------------------
New-ProductRow $ws 14 13 8 "TOBOLANZA 24 SOFT GELATIN CAPS." "1:0" "1" "49.00" "24.5000" "0:1"

# ---------------------------------------------------------------------------
# Step 2: insert CONTROLOC and ITOMASH before LIMITLESS (currently row 10).
# Clone formatting from row 11 (PHENADONE), the row immediately below.
# ---------------------------------------------------------------------------
New-ProductRow $ws 10 11 4 "CONTROLOC 40MG 14 GASTRORESISTANT TAB" "2:0" "1" "188.00" "188.0000" "1:0"
New-ProductRow $ws 11 12 5 "ITOMASH 50MG 30 TAB." "0:1" "1" "63.00" "20.7900" "0:1"

# ---------------------------------------------------------------------------
# Step 3: renumber the "م" sequence column for every product row (7..17) and
# fix up the row heights that changed as rows shuffled around.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11

$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75
$ws.Rows.Item(19).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 4: update the grand-total cell (old P15, now P18) and the generated
# timestamp in the footer (old A16, now A19).
# ---------------------------------------------------------------------------
$ws.Range("P18").Value = 511.62
$ws.Range("A19").Value = "Tuesday, 9 September, 2025 10:55 AM"
